$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object 'object[,]' 1,20
$row2[0,0] = "ECs"
$row2[0,1] = "Vcan"
$row2[0,2] = "Tlr2"
$row2[0,3] = "ECs"
$row2[0,4] = 3.0
$row2[0,5] = 1.0
$row2[0,6] = 7.172092999999999
$row2[0,7] = 21.516279
$row2[0,8] = 0.073573870768057
$row2[0,9] = 0.07357387076805699
$row2[0,10] = 3.0
$row2[0,11] = 1.0
$row2[0,12] = 57.24915866666667
$row2[0,13] = 171.747476
$row2[0,14] = 0.9704198736548433
$row2[0,15] = 0.9704198736548435
$row2[0,16] = 410.5962901290893
$row2[0,17] = 3695.366611161804
$row2[0,18] = 0.07139754637503565
$row2[0,19] = 0.07139754637503565
$ws.Range("A2:T2").Value = $row2

$row3 = New-Object 'object[,]' 1,20
$row3[0,0] = "ECs"
$row3[0,1] = "Vcan"
$row3[0,2] = "Tlr2"
$row3[0,3] = "FAPs"
$row3[0,4] = 3.0
$row3[0,5] = 1.0
$row3[0,6] = 7.172092999999999
$row3[0,7] = 21.516279
$row3[0,8] = 0.073573870768057
$row3[0,9] = 0.07357387076805699
$row3[0,10] = 3.0
$row3[0,11] = 1.0
$row3[0,12] = 1.569166
$row3[0,13] = 4.707498
$row3[0,14] = 0.02659864191768634
$row3[0,15] = 0.02659864191768634
$row3[0,16] = 11.254204484438
$row3[0,17] = 101.287840359942
$row3[0,18] = 0.001956965043057678
$row3[0,19] = 0.001956965043057678
$ws.Range("A3:T3").Value = $row3

$row4 = New-Object 'object[,]' 1,20
$row4[0,0] = "ECs"
$row4[0,1] = "Vcan"
$row4[0,2] = "Tlr2"
$row4[0,3] = "sCs"
$row4[0,4] = 3.0
$row4[0,5] = 1.0
$row4[0,6] = 7.172092999999999
$row4[0,7] = 21.516279
$row4[0,8] = 0.073573870768057
$row4[0,9] = 0.07357387076805699
$row4[0,10] = 2.0
$row4[0,11] = 0.6666666666666666
$row4[0,12] = 0.1758903333333333
$row4[0,13] = 0.527671
$row4[0,14] = 0.002981484427470275
$row4[0,15] = 0.002981484427470276
$row4[0,16] = 1.261501828467667
$row4[0,17] = 11.353516456209
$row4[0,18] = 0.0002193593499636725
$row4[0,19] = 0.0002193593499636725
$ws.Range("A4:T4").Value = $row4

$row5 = New-Object 'object[,]' 1,20
$row5[0,0] = "FAPs"
$row5[0,1] = "Vcan"
$row5[0,2] = "Tlr2"
$row5[0,3] = "ECs"
$row5[0,4] = 3.0
$row5[0,5] = 1.0
$row5[0,6] = 89.72947699999999
$row5[0,7] = 269.188431
$row5[0,8] = 0.9204767624852804
$row5[0,9] = 0.9204767624852804
$row5[0,10] = 3.0
$row5[0,11] = 1.0
$row5[0,12] = 57.24915866666667
$row5[0,13] = 171.747476
$row5[0,14] = 0.9704198736548433
$row5[0,15] = 0.9704198736548435
$row5[0,16] = 5136.937065850017
$row5[0,17] = 46232.43359265015
$row5[0,18] = 0.8932489435531851
$row5[0,19] = 0.8932489435531852
$ws.Range("A5:T5").Value = $row5

$row6 = New-Object 'object[,]' 1,20
$row6[0,0] = "FAPs"
$row6[0,1] = "Vcan"
$row6[0,2] = "Tlr2"
$row6[0,3] = "FAPs"
$row6[0,4] = 3.0
$row6[0,5] = 1.0
$row6[0,6] = 89.72947699999999
$row6[0,7] = 269.188431
$row6[0,8] = 0.9204767624852804
$row6[0,9] = 0.9204767624852804
$row6[0,10] = 3.0
$row6[0,11] = 1.0
$row6[0,12] = 1.569166
$row6[0,13] = 4.707498
$row6[0,14] = 0.02659864191768634
$row6[0,15] = 0.02659864191768634
$row6[0,16] = 140.800444506182
$row6[0,17] = 1267.204000555638
$row6[0,18] = 0.02448343179889719
$row6[0,19] = 0.02448343179889719
$ws.Range("A6:T6").Value = $row6

$row7 = New-Object 'object[,]' 1,20
$row7[0,0] = "FAPs"
$row7[0,1] = "Vcan"
$row7[0,2] = "Tlr2"
$row7[0,3] = "sCs"
$row7[0,4] = 3.0
$row7[0,5] = 1.0
$row7[0,6] = 89.72947699999999
$row7[0,7] = 269.188431
$row7[0,8] = 0.9204767624852804
$row7[0,9] = 0.9204767624852804
$row7[0,10] = 2.0
$row7[0,11] = 0.6666666666666666
$row7[0,12] = 0.1758903333333333
$row7[0,13] = 0.527671
$row7[0,14] = 0.002981484427470275
$row7[0,15] = 0.002981484427470276
$row7[0,16] = 15.78254761935567
$row7[0,17] = 142.042928574201
$row7[0,18] = 0.002744387133198119
$row7[0,19] = 0.002744387133198119
$ws.Range("A7:T7").Value = $row7

$row8 = New-Object 'object[,]' 1,20
$row8[0,0] = "sCs"
$row8[0,1] = "Vcan"
$row8[0,2] = "Tlr2"
$row8[0,3] = "ECs"
$row8[0,4] = 3.0
$row8[0,5] = 1.0
$row8[0,6] = 0.5799533333333334
$row8[0,7] = 1.73986
$row8[0,8] = 0.005949366746662454
$row8[0,9] = 0.005949366746662453
$row8[0,10] = 3.0
$row8[0,11] = 1.0
$row8[0,12] = 57.24915866666667
$row8[0,13] = 171.747476
$row8[0,14] = 0.9704198736548433
$row8[0,15] = 0.9704198736548435
$row8[0,16] = 33.20184039926223
$row8[0,17] = 298.8165635933601
$row8[0,18] = 0.005773383726622505
$row8[0,19] = 0.005773383726622505
$ws.Range("A8:T8").Value = $row8

$row9 = New-Object 'object[,]' 1,20
$row9[0,0] = "sCs"
$row9[0,1] = "Vcan"
$row9[0,2] = "Tlr2"
$row9[0,3] = "FAPs"
$row9[0,4] = 3.0
$row9[0,5] = 1.0
$row9[0,6] = 0.5799533333333334
$row9[0,7] = 1.73986
$row9[0,8] = 0.005949366746662454
$row9[0,9] = 0.005949366746662453
$row9[0,10] = 3.0
$row9[0,11] = 1.0
$row9[0,12] = 1.569166
$row9[0,13] = 4.707498
$row9[0,14] = 0.02659864191768634
$row9[0,15] = 0.02659864191768634
$row9[0,16] = 0.9100430522533335
$row9[0,17] = 8.190387470280001
$row9[0,18] = 0.0001582450757314652
$row9[0,19] = 0.0001582450757314651
$ws.Range("A9:T9").Value = $row9

$row10 = New-Object 'object[,]' 1,20
$row10[0,0] = "sCs"
$row10[0,1] = "Vcan"
$row10[0,2] = "Tlr2"
$row10[0,3] = "sCs"
$row10[0,4] = 3.0
$row10[0,5] = 1.0
$row10[0,6] = 0.5799533333333334
$row10[0,7] = 1.73986
$row10[0,8] = 0.005949366746662454
$row10[0,9] = 0.005949366746662453
$row10[0,10] = 2.0
$row10[0,11] = 0.6666666666666666
$row10[0,12] = 0.1758903333333333
$row10[0,13] = 0.527671
$row10[0,14] = 0.002981484427470275
$row10[0,15] = 0.002981484427470276
$row10[0,16] = 0.1020081851177778
$row10[0,17] = 0.9180736660600001
$row10[0,18] = 0.0000177379443084836
$row10[0,19] = 0.0000177379443084836
$ws.Range("A10:T10").Value = $row10
